$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of column A from "Savant Name" to "Baseball_Savant_Name"
$ws.Range("A1").Value = "Baseball_Savant_Name"

# Row 12 is an empty/stray data row (only has a Handedness value, no name/id) -
# delete it entirely so the rows below shift up by one.
$ws.Rows(12).Delete()

# Update the active selection to A2, matching the saved workbook view state.
$ws.Range("A2").Select()
